{"js": "// Update the date line and all the multiplication problems in the table\n// to match the new day's worksheet content.\nconst replacements = [\n  [\"2024-12-09 Monday\", \"2024-12-10 Tuesday\"],\n  [\"705\u00d78=\", \"548\u00d72=\"],\n  [\"524\u00d77=\", \"815\u00d76=\"],\n  [\"936\u00d74=\", \"983\u00d79=\"],\n  [\"392\u00d77=\", \"930\u00d76=\"],\n  [\"440\u00d73=\", \"307\u00d73=\"],\n  [\"944\u00d73=\", \"203\u00d78=\"],\n  [\"744\u00d78=\", \"261\u00d78=\"],\n  [\"534\u00d79=\", \"126\u00d77=\"],\n  [\"201\u00d78=\", \"982\u00d76=\"],\n  [\"608\u00d75=\", \"686\u00d72=\"],\n  [\"930\u00d74=\", \"279\u00d77=\"],\n  [\"618\u00d78=\", \"941\u00d76=\"],\n  [\"638\u00d78=\", \"376\u00d75=\"],\n  [\"556\u00d76=\", \"718\u00d79=\"],\n  [\"966\u00d77=\", \"932\u00d78=\"],\n  [\"715\u00d77=\", \"933\u00d78=\"],\n  [\"343\u00d74=\", \"403\u00d77=\"],\n  [\"878\u00d73=\", \"578\u00d79=\"],\n  [\"824\u00d72=\", \"788\u00d74=\"],\n  [\"713\u00d75=\", \"248\u00d75=\"],\n  [\"963\u00d79=\", \"314\u00d76=\"],\n  [\"802\u00d78=\", \"344\u00d79=\"],\n  [\"360\u00d75=\", \"200\u00d74=\"],\n  [\"189\u00d79=\", \"376\u00d73=\"],\n  [\"711\u00d72=\", \"312\u00d73=\"],\n];\n\nconst body = context.document.body;\n\nfor (const [oldText, newText] of replacements) {\n  const results = body.search(oldText, { matchCase: true, matchWholeWord: false });\n  results.load(\"items\");\n  await context.sync();\n\n  for (const item of results.items) {\n    item.insertText(newText, Word.InsertLocation.replace);\n  }\n  await context.sync();\n}\n", "ps1": "# Update the date line and all the multiplication problems in the table\n# to match the new day's worksheet content.\n$d = $word.ActiveDocument\n\n$replacements = @(\n    @(\"2024-12-09 Monday\", \"2024-12-10 Tuesday\"),\n    @(\"705\u00d78=\", \"548\u00d72=\"),\n    @(\"524\u00d77=\", \"815\u00d76=\"),\n    @(\"936\u00d74=\", \"983\u00d79=\"),\n    @(\"392\u00d77=\", \"930\u00d76=\"),\n    @(\"440\u00d73=\", \"307\u00d73=\"),\n    @(\"944\u00d73=\", \"203\u00d78=\"),\n    @(\"744\u00d78=\", \"261\u00d78=\"),\n    @(\"534\u00d79=\", \"126\u00d77=\"),\n    @(\"201\u00d78=\", \"982\u00d76=\"),\n    @(\"608\u00d75=\", \"686\u00d72=\"),\n    @(\"930\u00d74=\", \"279\u00d77=\"),\n    @(\"618\u00d78=\", \"941\u00d76=\"),\n    @(\"638\u00d78=\", \"376\u00d75=\"),\n    @(\"556\u00d76=\", \"718\u00d79=\"),\n    @(\"966\u00d77=\", \"932\u00d78=\"),\n    @(\"715\u00d77=\", \"933\u00d78=\"),\n    @(\"343\u00d74=\", \"403\u00d77=\"),\n    @(\"878\u00d73=\", \"578\u00d79=\"),\n    @(\"824\u00d72=\", \"788\u00d74=\"),\n    @(\"713\u00d75=\", \"248\u00d75=\"),\n    @(\"963\u00d79=\", \"314\u00d76=\"),\n    @(\"802\u00d78=\", \"344\u00d79=\"),\n    @(\"360\u00d75=\", \"200\u00d74=\"),\n    @(\"189\u00d79=\", \"376\u00d73=\"),\n    @(\"711\u00d72=\", \"312\u00d73=\")\n)\n\nforeach ($pair in $replacements) {\n    $old = $pair[0]\n    $new = $pair[1]\n    $find = $d.Content.Find\n    $find.ClearFormatting()\n    $find.Replacement.ClearFormatting()\n    $find.Text = $old\n    $find.Replacement.Text = $new\n    $find.Execute($old, $false, $false, $false, $false, $false, $true, 1, $false, $new, 2)\n}\n"}
